$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "M1"
$ws.Range("B2").Value = "Areg"
$ws.Range("C2").Value = "Erbb3"
$ws.Range("D2").Value = "FAPs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 1.484834666666667
$ws.Range("H2").Value = 4.454504
$ws.Range("I2").Value = 0.4307162850350085
$ws.Range("J2").Value = 0.4307162850350084
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("M2").Value = 0.2048153333333333
$ws.Range("N2").Value = 0.614446
$ws.Range("O2").Value = 0.06238051989214988
$ws.Range("P2").Value = 0.06238051989214988
$ws.Range("Q2").Value = 0.3041169071982223
$ws.Range("R2").Value = 2.737052164784
$ws.Range("S2").Value = 0.02686830578649925
$ws.Range("T2").Value = 0.02686830578649924

$ws.Range("A3").Value = "M1"
$ws.Range("B3").Value = "Areg"
$ws.Range("C3").Value = "Erbb3"
$ws.Range("D3").Value = "M1"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 1.484834666666667
$ws.Range("H3").Value = 4.454504
$ws.Range("I3").Value = 0.4307162850350085
$ws.Range("J3").Value = 0.4307162850350084
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 0.7041406666666666
$ws.Range("N3").Value = 2.112422
$ws.Range("O3").Value = 0.2144598265618379
$ws.Range("P3").Value = 0.2144598265618379
$ws.Range("Q3").Value = 1.045532472076445
$ws.Range("R3").Value = 9.409792248688
$ws.Range("S3").Value = 0.09237133978596704
$ws.Range("T3").Value = 0.09237133978596702

$ws.Range("A4").Value = "M1"
$ws.Range("B4").Value = "Areg"
$ws.Range("C4").Value = "Erbb3"
$ws.Range("D4").Value = "M2"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 1.484834666666667
$ws.Range("H4").Value = 4.454504
$ws.Range("I4").Value = 0.4307162850350085
$ws.Range("J4").Value = 0.4307162850350084
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 0.7300083333333333
$ws.Range("N4").Value = 2.190025
$ws.Range("O4").Value = 0.2223383309140356
$ws.Range("P4").Value = 0.2223383309140356
$ws.Range("Q4").Value = 1.083941680288889
$ws.Range("R4").Value = 9.7554751226
$ws.Range("S4").Value = 0.09576473991217781
$ws.Range("T4").Value = 0.0957647399121778

$ws.Range("A5").Value = "M1"
$ws.Range("B5").Value = "Areg"
$ws.Range("C5").Value = "Erbb3"
$ws.Range("D5").Value = "sCs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 1.484834666666667
$ws.Range("H5").Value = 4.454504
$ws.Range("I5").Value = 0.4307162850350085
$ws.Range("J5").Value = 0.4307162850350084
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 1.644357666666667
$ws.Range("N5").Value = 4.933073
$ws.Range("O5").Value = 0.5008213226319767
$ws.Range("P5").Value = 0.5008213226319767
$ws.Range("Q5").Value = 2.441599267865778
$ws.Range("R5").Value = 21.974393410792
$ws.Range("S5").Value = 0.2157118995503644
$ws.Range("T5").Value = 0.2157118995503644

$ws.Range("A6").Value = "M2"
$ws.Range("B6").Value = "Areg"
$ws.Range("C6").Value = "Erbb3"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 1.146843333333333
$ws.Range("H6").Value = 3.44053
$ws.Range("I6").Value = 0.332672795927784
$ws.Range("J6").Value = 0.332672795927784
$ws.Range("K6").Value = 2
$ws.Range("L6").Value = 0.6666666666666666
$ws.Range("M6").Value = 0.2048153333333333
$ws.Range("N6").Value = 0.614446
$ws.Range("O6").Value = 0.06238051989214988
$ws.Range("P6").Value = 0.06238051989214988
$ws.Range("Q6").Value = 0.2348910995977778
$ws.Range("R6").Value = 2.11401989638
$ws.Range("S6").Value = 0.02075230196395025
$ws.Range("T6").Value = 0.02075230196395025

$ws.Range("A7").Value = "M2"
$ws.Range("B7").Value = "Areg"
$ws.Range("C7").Value = "Erbb3"
$ws.Range("D7").Value = "M1"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 1.146843333333333
$ws.Range("H7").Value = 3.44053
$ws.Range("I7").Value = 0.332672795927784
$ws.Range("J7").Value = 0.332672795927784
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 0.7041406666666666
$ws.Range("N7").Value = 2.112422
$ws.Range("O7").Value = 0.2144598265618379
$ws.Range("P7").Value = 0.2144598265618379
$ws.Range("Q7").Value = 0.8075390292955554
$ws.Range("R7").Value = 7.26785126366
$ws.Range("S7").Value = 0.07134495011651423
$ws.Range("T7").Value = 0.07134495011651423

$ws.Range("A8").Value = "M2"
$ws.Range("B8").Value = "Areg"
$ws.Range("C8").Value = "Erbb3"
$ws.Range("D8").Value = "M2"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 1.146843333333333
$ws.Range("H8").Value = 3.44053
$ws.Range("I8").Value = 0.332672795927784
$ws.Range("J8").Value = 0.332672795927784
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 0.7300083333333333
$ws.Range("N8").Value = 2.190025
$ws.Range("O8").Value = 0.2223383309140356
$ws.Range("P8").Value = 0.2223383309140356
$ws.Range("Q8").Value = 0.837205190361111
$ws.Range("R8").Value = 7.534846713249999
$ws.Range("S8").Value = 0.07396591418708907
$ws.Range("T8").Value = 0.07396591418708907

$ws.Range("A9").Value = "M2"
$ws.Range("B9").Value = "Areg"
$ws.Range("C9").Value = "Erbb3"
$ws.Range("D9").Value = "sCs"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 1.146843333333333
$ws.Range("H9").Value = 3.44053
$ws.Range("I9").Value = 0.332672795927784
$ws.Range("J9").Value = 0.332672795927784
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 1.644357666666667
$ws.Range("N9").Value = 4.933073
$ws.Range("O9").Value = 0.5008213226319767
$ws.Range("P9").Value = 0.5008213226319767
$ws.Range("Q9").Value = 1.885820627632222
$ws.Range("R9").Value = 16.97238564869
$ws.Range("S9").Value = 0.1666096296602304
$ws.Range("T9").Value = 0.1666096296602304

$ws.Range("A10").Value = "sCs"
$ws.Range("B10").Value = "Areg"
$ws.Range("C10").Value = "Erbb3"
$ws.Range("D10").Value = "FAPs"
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 0.8156833333333333
$ws.Range("H10").Value = 2.44705
$ws.Range("I10").Value = 0.2366109190372076
$ws.Range("J10").Value = 0.2366109190372076
$ws.Range("K10").Value = 2
$ws.Range("L10").Value = 0.6666666666666666
$ws.Range("M10").Value = 0.2048153333333333
$ws.Range("N10").Value = 0.614446
$ws.Range("O10").Value = 0.06238051989214988
$ws.Range("P10").Value = 0.06238051989214988
$ws.Range("Q10").Value = 0.1670644538111111
$ws.Range("R10").Value = 1.5035800843
$ws.Range("S10").Value = 0.01475991214170039
$ws.Range("T10").Value = 0.01475991214170039

$ws.Range("A11").Value = "sCs"
$ws.Range("B11").Value = "Areg"
$ws.Range("C11").Value = "Erbb3"
$ws.Range("D11").Value = "M1"
$ws.Range("E11").Value = 3
$ws.Range("F11").Value = 1
$ws.Range("G11").Value = 0.8156833333333333
$ws.Range("H11").Value = 2.44705
$ws.Range("I11").Value = 0.2366109190372076
$ws.Range("J11").Value = 0.2366109190372076
$ws.Range("K11").Value = 3
$ws.Range("L11").Value = 1
$ws.Range("M11").Value = 0.7041406666666666
$ws.Range("N11").Value = 2.112422
$ws.Range("O11").Value = 0.2144598265618379
$ws.Range("P11").Value = 0.2144598265618379
$ws.Range("Q11").Value = 0.5743558061222221
$ws.Range("R11").Value = 5.1692022551
$ws.Range("S11").Value = 0.0507435366593566
$ws.Range("T11").Value = 0.0507435366593566

$ws.Range("A12").Value = "sCs"
$ws.Range("B12").Value = "Areg"
$ws.Range("C12").Value = "Erbb3"
$ws.Range("D12").Value = "M2"
$ws.Range("E12").Value = 3
$ws.Range("F12").Value = 1
$ws.Range("G12").Value = 0.8156833333333333
$ws.Range("H12").Value = 2.44705
$ws.Range("I12").Value = 0.2366109190372076
$ws.Range("J12").Value = 0.2366109190372076
$ws.Range("K12").Value = 3
$ws.Range("L12").Value = 1
$ws.Range("M12").Value = 0.7300083333333333
$ws.Range("N12").Value = 2.190025
$ws.Range("O12").Value = 0.2223383309140356
$ws.Range("P12").Value = 0.2223383309140356
$ws.Range("Q12").Value = 0.5954556306944444
$ws.Range("R12").Value = 5.35910067625
$ws.Range("S12").Value = 0.05260767681476875
$ws.Range("T12").Value = 0.05260767681476875

$ws.Range("A13").Value = "sCs"
$ws.Range("B13").Value = "Areg"
$ws.Range("C13").Value = "Erbb3"
$ws.Range("D13").Value = "sCs"
$ws.Range("E13").Value = 3
$ws.Range("F13").Value = 1
$ws.Range("G13").Value = 0.8156833333333333
$ws.Range("H13").Value = 2.44705
$ws.Range("I13").Value = 0.2366109190372076
$ws.Range("J13").Value = 0.2366109190372076
$ws.Range("K13").Value = 3
$ws.Range("L13").Value = 1
$ws.Range("M13").Value = 1.644357666666667
$ws.Range("N13").Value = 4.933073
$ws.Range("O13").Value = 0.5008213226319767
$ws.Range("P13").Value = 0.5008213226319767
$ws.Range("Q13").Value = 1.341275142738889
$ws.Range("R13").Value = 12.07147628465
$ws.Range("S13").Value = 0.1184997934213818
$ws.Range("T13").Value = 0.1184997934213818

